$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '61.347.13'
$ws.Range("E2").Value = '  +0.96%  '
$ws.Range("D3").Value = '3.429.08'
$ws.Range("E3").Value = '  +1.71%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '575.66'
$ws.Range("E5").Value = '  +1.09%  '
$ws.Range("D6").Value = '145.27'
$ws.Range("E6").Value = '  +7.03%  '
$ws.Range("D7").Value = '3.428.98'
$ws.Range("E7").Value = '  +1.76%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  +2.33%  '
$ws.Range("E10").Value = '  +0.99%  '
$ws.Range("E11").Value = '  +3.47%  '
$ws.Range("E12").Value = '  +2.09%  '
$ws.Range("D13").Value = '4.022.87'
$ws.Range("E13").Value = '  +1.96%  '
$ws.Range("D14").Value = '27.92'
$ws.Range("E14").Value = '  +7.55%  '
$ws.Range("E15").Value = '  -0.60%  '
$ws.Range("E16").Value = '  +1.91%  '
$ws.Range("D17").Value = '3.429.65'
$ws.Range("E17").Value = '  +1.73%  '
$ws.Range("D18").Value = '61.472.42'
$ws.Range("E18").Value = '  +1.07%  '
$ws.Range("D19").Value = '6.29'
$ws.Range("E19").Value = '  +8.32%  '
$ws.Range("D20").Value = '14.18'
$ws.Range("E20").Value = '  +3.18%  '
$ws.Range("D21").Value = '9.40'
$ws.Range("E21").Value = '  +2.27%  '
$ws.Range("D22").Value = '395.06'
$ws.Range("E22").Value = '  +6.32%  '
$ws.Range("D23").Value = '0.565'
$ws.Range("E23").Value = '  +3.30%  '
$ws.Range("D24").Value = '73.34'
$ws.Range("E24").Value = '  +3.79%  '
$ws.Range("E25").Value = '  -0.31%  '
$ws.Range("E26").Value = '  +0.44%  '
$ws.Range("E27").Value = '  +0.28%  '
$ws.Range("D28").Value = '3.571.88'
$ws.Range("E28").Value = '  +1.79%  '
$ws.Range("D29").Value = '0.180'
$ws.Range("E29").Value = '  +3.10%  '
$ws.Range("E30").Value = '  +3.52%  '
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("E32").Value = '  +2.32%  '
$ws.Range("E33").Value = '  -6.68%  '
$ws.Range("E34").Value = '  +2.50%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").Value = '23.95'
$ws.Range("E36").Value = '  +2.96%  '
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").Value = '7.02'
$ws.Range("E37").Value = '  +3.99%  '
$ws.Range("B38").Value = 'RenzoRestakedETH'
$ws.Range("C38").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D38").Value = '3.462.03'
$ws.Range("E38").Value = '  +2.15%  '
$ws.Range("E39").Value = '  +1.74%  '
$ws.Range("D40").Value = '5.11'
$ws.Range("E40").Value = '  +0.47%  '
$ws.Range("D41").Value = '167.74'
$ws.Range("E41").Value = '  +2.04%  '
$ws.Range("D42").Value = '0.0783'
$ws.Range("E42").Value = '  +3.11%  '
$ws.Range("D43").Value = '26.86'
$ws.Range("E43").Value = '  +5.54%  '
$ws.Range("D44").Value = '0.798'
$ws.Range("E44").Value = '  +3.75%  '
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("E46").Value = '  +0.09%  '
$ws.Range("D47").Value = '4.48'
$ws.Range("E47").Value = '  +3.78%  '
$ws.Range("E48").Value = '  +0.65%  '
$ws.Range("D49").Value = '2.590.96'
$ws.Range("E49").Value = '  +3.13%  '
$ws.Range("E50").Value = '  -0.10%  '
$ws.Range("E51").Value = '  +2.82%  '
